$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the content of row 2 and row 3 (columns A through E), which represent
# two news-article records in the historical distance time-bucket dataset.
# After the swap, the "Naroda Patiya: Maya Kodnani jailed..." record (originally
# row 3) becomes row 2, and the "Arjun Modhwadia on Twitter" record (originally
# row 2) becomes row 3 - matching the updated historical-distance values
# (3837 / 3835) used for the time bucket analysis.

$cols = @("A", "B", "C", "D", "E")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell3 = $ws.Range($col + "3")

    $v2 = $cell2.Value2
    $v3 = $cell3.Value2

    $cell2.Value2 = $v3
    $cell3.Value2 = $v2
}
